{"js": "// Update the multiplication problems in the table to a new set of\n// operands, per the commit's regenerated \"output\" table.\n// Each left-hand text is unique within the document, so a simple\n// search-and-replace (matchCase, whole match) per pair is safe and\n// avoids any row/column index bookkeeping.\n\nconst replacements = [\n  [\"52\u00d719=\", \"35\u00d728=\"],\n  [\"33\u00d791=\", \"91\u00d741=\"],\n  [\"86\u00d762=\", \"48\u00d778=\"],\n  [\"57\u00d759=\", \"17\u00d734=\"],\n  [\"50\u00d785=\", \"43\u00d715=\"],\n  [\"17\u00d730=\", \"96\u00d782=\"],\n  [\"12\u00d788=\", \"84\u00d784=\"],\n  [\"93\u00d713=\", \"97\u00d782=\"],\n  [\"15\u00d739=\", \"88\u00d776=\"],\n  [\"20\u00d726=\", \"84\u00d728=\"],\n  [\"31\u00d798=\", \"37\u00d739=\"],\n  [\"54\u00d756=\", \"88\u00d795=\"],\n  [\"93\u00d735=\", \"22\u00d778=\"],\n  [\"80\u00d711=\", \"57\u00d762=\"],\n  [\"84\u00d734=\", \"30\u00d792=\"],\n  [\"27\u00d723=\", \"16\u00d722=\"],\n  [\"47\u00d798=\", \"67\u00d774=\"],\n  [\"15\u00d795=\", \"82\u00d733=\"],\n  [\"62\u00d753=\", \"55\u00d738=\"],\n  [\"63\u00d759=\", \"67\u00d738=\"],\n  [\"37\u00d759=\", \"13\u00d797=\"],\n  [\"26\u00d777=\", \"22\u00d715=\"],\n  [\"41\u00d740=\", \"18\u00d789=\"],\n  [\"25\u00d773=\", \"97\u00d748=\"],\n  [\"55\u00d777=\", \"20\u00d769=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const found = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  await context.sync();\n\n  if (found.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  for (const range of found.items) {\n    range.insertText(newText, \"Replace\");\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the multiplication problems in the table to a new set of\n# operands, per the commit's regenerated \"output\" table.\n# Each left-hand text is unique within the document, so Find/Replace\n# (one pair at a time, matching the whole \"NN\u00d7NN=\" token) is safe and\n# avoids any row/column index bookkeeping.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"52\u00d719=\", \"35\u00d728=\"),\n    @(\"33\u00d791=\", \"91\u00d741=\"),\n    @(\"86\u00d762=\", \"48\u00d778=\"),\n    @(\"57\u00d759=\", \"17\u00d734=\"),\n    @(\"50\u00d785=\", \"43\u00d715=\"),\n    @(\"17\u00d730=\", \"96\u00d782=\"),\n    @(\"12\u00d788=\", \"84\u00d784=\"),\n    @(\"93\u00d713=\", \"97\u00d782=\"),\n    @(\"15\u00d739=\", \"88\u00d776=\"),\n    @(\"20\u00d726=\", \"84\u00d728=\"),\n    @(\"31\u00d798=\", \"37\u00d739=\"),\n    @(\"54\u00d756=\", \"88\u00d795=\"),\n    @(\"93\u00d735=\", \"22\u00d778=\"),\n    @(\"80\u00d711=\", \"57\u00d762=\"),\n    @(\"84\u00d734=\", \"30\u00d792=\"),\n    @(\"27\u00d723=\", \"16\u00d722=\"),\n    @(\"47\u00d798=\", \"67\u00d774=\"),\n    @(\"15\u00d795=\", \"82\u00d733=\"),\n    @(\"62\u00d753=\", \"55\u00d738=\"),\n    @(\"63\u00d759=\", \"67\u00d738=\"),\n    @(\"37\u00d759=\", \"13\u00d797=\"),\n    @(\"26\u00d777=\", \"22\u00d715=\"),\n    @(\"41\u00d740=\", \"18\u00d789=\"),\n    @(\"25\u00d773=\", \"97\u00d748=\"),\n    @(\"55\u00d777=\", \"20\u00d769=\")\n)\n\nforeach ($pair in $replacements) {\n    $old = $pair[0]\n    $new = $pair[1]\n\n    $rng = $d.Content\n    $find = $rng.Find\n    $find.ClearFormatting()\n    $find.Text = $old\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $new\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Forward = $true\n    $find.Wrap = 1\n\n    $find.Execute(\n        $old,\n        $true,\n        $false,\n        $false,\n        $false,\n        $false,\n        $true,\n        1,\n        $false,\n        $new,\n        2\n    )\n}\n\n$d.Save()\n"}
